$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Stellen Sie sich vor, Sie wären ein Sicherheitsperson der auffällige Aktivitäten in einem Unternehmen überwacht. Ihre Aufgabe erfordert ständige Aufmerksamkeit und schnelle Reaktionen, wenn etwas Verdächtiges passiert.   `n`nWir untersuchen in unserem Labor Aufmerksamkeit und Reaktionsgeschwindigkeit und in diesem Experiment bitten wir Sie, die Rolle des Sicherheitsmannes zu spielen.  `n`nGenauer werden Sie eine Reihe von Dingen auf dem Computerbildschirm beobachten und so schnell wie möglich reagieren, indem Sie die Leertaste drücken, wenn ein Zielgegenstand auftaucht.  `n`n`nDrücken Sie die Leertaste, um fortzufahren."

$ws.Range("A9").Value = "Als Nächstes werden Ihnen 30 Paare aus Ziel- und Füll-Wesen aus der Überwachungsaufgabe gezeigt und wir bitten Sie anzugeben, welches Sie lieber mögen.  `n`nSie brauchen keinen Grund, um eines lieber als das andere zu mögen. Teilen Sie uns einfach mit, was Ihr Bauchgefühl ist.`n`nUns interessiert, ob die Angenehmheit oder Unangenehmheit der Wesen die Fähigkeit beeinflusst, sie aufmerksam zu beobachten und schnell auf sie zu reagieren. Daher benötigen wir Ihre Angabe, welches Sie lieber mögen.  `n`nNicht vergessen: Sie brauchen keinen Grund, um eines lieber als das andere zu mögen, also folgen Sie einfach Ihrem Bauchgefühl. Bitte antworten Sie zügig.`n`n`nDrücken Sie die Leertaste, um fortzufahren. "
